$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.954.69'
$ws.Range('E2').Value = '  +0.40%  '
$ws.Range('D3').Value = '1.818.03'
$ws.Range('E3').Value = '  +0.49%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '310.02'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.09%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.002'
$ws.Range('D6').ClearFormats()
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4697'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.78%  '
$ws.Range('E8').Value = '  -0.81%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07353'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.09%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.8728'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.02%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '20.28'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.79%  '
$ws.Range('D12').Value = '1.815.25'
$ws.Range('E12').Value = '  -1.44%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.408'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.99%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.07113'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.58%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.510'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.07%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '91.40'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.07%  '
$ws.Range('E17').Value = '  +0.09%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.000008705'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.02%  '
$ws.Range('E19').Value = '  +0.09%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '14.66'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.20%  '
$ws.Range('D21').Value = '26.968.82'
$ws.Range('E21').Value = '  +0.37%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.294'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.40%  '
$ws.Range('E23').Value = '  +0.81%  '
$ws.Range('D24').Value = '2.042.28'
$ws.Range('E24').Value = '  -0.71%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.894'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.58%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '150.83'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.52%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '18.42'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.34%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.149'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.25%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '5.259'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.83%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '117.00'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.98%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.08902'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.04%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.7604'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.05%  '
$ws.Range('E33').Value = '  +0.75%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.502'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +1.11%  '
$ws.Range('E35').Value = '  -0.30%  '
$ws.Range('E36').Value = '  +0.11%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.094'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.38%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.05293'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +0.81%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.01946'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.81%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.969'
$ws.Range('D40').ClearFormats()
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.398'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -1.33%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.5290'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.30%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '7.151'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.11%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.1654'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.62%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '8.444'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.06%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.4870'
$ws.Range('D46').ClearFormats()
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '10.43'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +1.05%  '
$ws.Range('E48').Value = '  +0.10%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '103.47'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.44%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.664'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.45%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.06296'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.29%  '
